# edit.ps1 - applies the "week 11 / week 10 assignment update" edit
# described by the diff to before.pptx via PowerPoint COM interop.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 2 ("Housekeeping"): extend the "Ideas? " bullet with more text.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$content2 = $s2.Shapes.Item(2)
$tr2 = $content2.TextFrame.TextRange
$ideasPara = $tr2.Paragraphs(4)
$ideasPara.Text = "Ideas? Bayesian methods, multivariate methods, networks ("
$r = $ideasPara.InsertAfter("prob")
$r = $r.InsertAfter(" not), ")

# ---------------------------------------------------------------------
# 2) Restructure the "Problems" / "Advanced topics" slides (currently
#    slides 31 and 32) into three slides:
#      31: Advanced topics (unchanged content, duplicated)
#      32: Assignment (brand new GLM assignment content)
#      33: Problems (unchanged content, the original "Problems" slide)
# ---------------------------------------------------------------------

# Duplicate "Advanced topics" (currently slide 32); the copy lands at 33.
$dup = $p.Slides.Item(32).Duplicate()

# Move the original "Problems" slide (currently slide 31) to the end
# of this trio (position 33); "Advanced topics" (orig) shifts to 31,
# the "Advanced topics" duplicate shifts to 32.
$p.Slides.Item(31).MoveTo(33)

# Now:
#   slide 31 = Advanced topics (original, untouched) -> matches target
#   slide 32 = Advanced topics (duplicate) -> turn into "Assignment"
#   slide 33 = Problems (original, untouched) -> matches target

# --- Edit slide 32 into the new "Assignment" slide -----------------
$sAssign = $p.Slides.Item(32)

$title = $sAssign.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Assignment"

$body = $sAssign.Shapes.Item(2)

# Reposition / resize the content placeholder.
$body.Left = 66.0
$body.Top = 112.18606299212598
$body.Width = 828.0
$body.Height = 374.18898017795277

$btr = $body.TextFrame.TextRange
$btr.Text = "Part 1 - GLMs`rMake a generalized linear model (preferably with more than one variable) for one of your hypotheses. Articulate which hypothesis you are testing.`rExplain what the R output is telling you about your data, in relation to your hypothesis.`rPlot your model (e.g. using predict) and overlay the model on top of the underlying data. Remember that you will need to use " + [char]8220 + "type=response" + [char]8221 + ".`rWrite a results statement (as you would in a scientific paper). If you need to reference a statistical table, you can include this result statement and table as a separate document that you upload with your text entry. (You can title this LASTNAME_WEEK10_Results.) Please remember to also commit and push to github. `rYou will turn in this assignment in two weeks with model comparisons`r`r`r"

# Paragraph 1: "Part 1 - GLMs" - no bullet, flush-left, underlined.
$p1 = $btr.Paragraphs(1)
$p1.ParagraphFormat.Bullet.Type = 0
$p1.Font.Underline = -1

Write-Output "done"
